$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.475.44"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "2.339.12"
$ws.Range("E3").Value = "  -1.75%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.57"
$ws.Range("E5").Value = "  -4.62%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "107.66"
$ws.Range("E6").Value = "  +1.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.631"
$ws.Range("E7").Value = "  -1.84%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -6.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.88"
$ws.Range("E10").Value = "  -1.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0923"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.46"
$ws.Range("E12").Value = "  -2.06%  "
$ws.Range("E13").Value = "  -0.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.988"
$ws.Range("E14").Value = "  -5.65%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.85"
$ws.Range("E15").Value = "  -7.59%  "
$ws.Range("D16").Value = "2.692.19"
$ws.Range("E16").Value = "  -1.86%  "
$ws.Range("D17").Value = "2.336.93"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("D18").Value = "42.393.19"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.59"
$ws.Range("E19").Value = "  -4.58%  "
$ws.Range("E20").Value = "  -3.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "75.90"
$ws.Range("E21").Value = "  -1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.57"
$ws.Range("E22").Value = "  -2.38%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "255.91"
$ws.Range("E23").Value = "  -8.29%  "
$ws.Range("E24").Value = "  -5.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.34"
$ws.Range("E25").Value = "  -3.09%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.11%  "
$ws.Range("E27").Value = "  -3.85%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "22.78"
$ws.Range("E28").Value = "  -1.80%  "
$ws.Range("E29").Value = "  +1.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "173.54"
$ws.Range("E30").Value = "  -0.68%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "36.49"
$ws.Range("E31").Value = "  -4.10%  "
$ws.Range("E32").Value = "  -4.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "6.03"
$ws.Range("E33").Value = "  +2.58%  "
$ws.Range("E34").Value = "  -10.15%  "
$ws.Range("E35").Value = "  +17.15%  "
$ws.Range("E36").Value = "  -1.98%  "
$ws.Range("E37").Value = "  -5.83%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0361"
$ws.Range("E38").Value = "  -2.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.88"
$ws.Range("E39").Value = "  -9.53%  "
$ws.Range("E40").Value = "  -6.60%  "
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.46"
$ws.Range("E42").Value = "  -7.69%  "
$ws.Range("B43").Value = "MultiversX"
$ws.Range("C43").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.39"
$ws.Range("E43").Value = "  +0.52%  "
$ws.Range("E44").Value = "  -0.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "11.92"
$ws.Range("E45").Value = "  -4.70%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "111.23"
$ws.Range("E46").Value = "  -9.94%  "
$ws.Range("E47").Value = "  -2.84%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.29"
$ws.Range("E48").Value = "  -10.12%  "
$ws.Range("E49").Value = "  -1.97%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "74.06"
$ws.Range("E50").Value = "  +1.83%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.28"
$ws.Range("E51").Value = "  -2.79%  "
